# Atualização de bases das ligas, do dia: 11-04-2024 às 00:31
#
# The source rows were re-ordered (their business-key / "id" column B and
# everything through column AC got shuffled among a handful of rows, while
# the running index in column A stayed put). This script reproduces that
# reshuffle by rotating the B:AC payload among the affected rows, working
# from an in-memory snapshot so the rotations (including the two 4-row
# cycles) apply atomically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is an ordered cycle of row numbers: the new B:AC content
# of cycle[i] becomes the old B:AC content of cycle[i+1] (wrapping around).
$cycles = @(
    ,@(23, 24)
    ,@(26, 27)
    ,@(112, 113)
    ,@(118, 120, 119, 121)
    ,@(125, 128, 126, 129)
    ,@(138, 139)
    ,@(156, 157)
    ,@(194, 195)
)

foreach ($cycle in $cycles) {
    $n = $cycle.Length

    # Snapshot every row's current B:AC payload before writing anything back.
    $stored = @{}
    foreach ($r in $cycle) {
        $rng = $ws.Range("B$r`:AC$r")
        $stored[$r] = $rng.Value2
    }

    for ($i = 0; $i -lt $n; $i++) {
        $target = $cycle[$i]
        $src = $cycle[(($i + 1) % $n)]
        $rng = $ws.Range("B$target`:AC$target")
        $rng.Value2 = $stored[$src]
    }
}
